$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph near the top of the document
#    (the paragraph right after the H1 title: "Meta description: Read our
#    review of Bingo Billions...").
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph ("Play Bingo Billions Free: Fun Gameplay
#    and Great Winning Potential") right before the final "Prompt:" paragraph,
#    built via literal OOXML so the run/paragraph formatting matches exactly
#    (no inherited italics, no stray pPr/style overrides).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$insertPoint = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)

$newParaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Bingo Billions Free: Fun Gameplay and Great Winning Potential</w:t></w:r></w:p>"
$insertPoint.InsertXML($newParaXml)

$headingText = "Play Bingo Billions Free: Fun Gameplay and Great Winning Potential"
$splitPos = $lastPara.Range.Start + $headingText.Length
$splitRange = $d.Range($splitPos, $splitPos)
$splitRange.InsertParagraphAfter()

# 3. Replace the old "Prompt: ..." text (now in the following paragraph)
#    with the meta-description copy.
$d.Content.Find.Execute(
    "Prompt: Create a feature image for Bingo Billions that reflects its theme and style. The image should be in a cartoon style and feature a happy Maya warrior with glasses, highlighting the fun and entertaining nature of the game while tying in the bingo concept. The image can include elements such as gold coins, bingo balls, and slot reels to represent the game's features and symbols. The overall tone should be upbeat and vibrant, showcasing the colorful and energetic design of Bingo Billions while catching the attention of potential players.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our review of Bingo Billions, a bingo-themed slot game with fun gameplay and great winning potential. Play free and enjoy frequent payouts and free spins.",
    2)
